$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.892.30'
$ws.Range("E2").Value = '  +4.70%  '

$ws.Range("D3").Value = '3.354.16'
$ws.Range("E3").Value = '  +4.98%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = "'557.62"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.64%  '

$ws.Range("D6").Value = "'153.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.73%  '

$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  +0.93%  '

$ws.Range("D9").Value = "'7.52"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("E10").Value = '  +4.30%  '

$ws.Range("E11").Value = '  +1.89%  '

$ws.Range("D12").Value = '3.932.22'
$ws.Range("E12").Value = '  +4.98%  '

$ws.Range("E13").Value = '  +0.62%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").Value = "'27.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.02%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = "'0.0000181"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.60%  '

$ws.Range("D16").Value = '62.971.35'
$ws.Range("E16").Value = '  +4.78%  '

$ws.Range("D17").Value = '3.367.48'
$ws.Range("E17").Value = '  +5.30%  '

$ws.Range("D18").Value = "'6.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.76%  '

$ws.Range("D19").Value = "'13.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.59%  '

$ws.Range("E20").Value = '  +1.03%  '

$ws.Range("D21").Value = "'388.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.90%  '

$ws.Range("D22").Value = "'0.542"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.42%  '

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.05%  '

$ws.Range("D24").Value = "'70.57"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").Value = "'0.179"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.82%  '

$ws.Range("D26").Value = "'8.90"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.36%  '

$ws.Range("E27").Value = '  +7.60%  '

$ws.Range("E28").Value = '  +0.18%  '

$ws.Range("D29").Value = "'6.63"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +7.29%  '

$ws.Range("E30").Value = '  +4.50%  '

$ws.Range("E31").Value = '  +5.18%  '

$ws.Range("D32").Value = "'23.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.01%  '

$ws.Range("D33").Value = "'1.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.30%  '

$ws.Range("E34").Value = '  +2.10%  '

$ws.Range("D35").Value = "'160.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.76%  '

$ws.Range("E36").Value = '  +9.45%  '

$ws.Range("E37").Value = '  +11.96%  '

$ws.Range("D38").Value = "'27.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.15%  '

$ws.Range("E39").Value = '  +3.81%  '

$ws.Range("D40").Value = '2.826.68'
$ws.Range("E40").Value = '  +1.84%  '

$ws.Range("D41").Value = "'0.0310"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.26%  '

$ws.Range("E42").Value = '  +1.59%  '

$ws.Range("D43").Value = "'40.78"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.56%  '

$ws.Range("D44").Value = "'0.747"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.38%  '

$ws.Range("E45").Value = '  +4.70%  '

$ws.Range("D46").Value = "'22.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +7.89%  '

$ws.Range("D47").Value = '3.398.89'
$ws.Range("E47").Value = '  +5.09%  '

$ws.Range("E48").Value = '  +2.40%  '

$ws.Range("E49").Value = '  +2.07%  '

$ws.Range("E50").Value = '  +0.55%  '

$ws.Range("D51").Value = "'278.96"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +5.84%  '
